$wb = $excel.ActiveWorkbook

# Sheet order (1-based, matches xl/workbook.xml):
#   1 Funciones_Objetivo
#   2 Restricciones_del_lider
#   3 Restricciones_del_follower
#   4 Punto_modificado
#   5 Vector_bf
#   6 Vector_BF
#   7 Vector_Alpha
# Using numeric indices avoids ambiguity between the case-insensitive
# sheet names "Vector_bf" and "Vector_BF".

# NOTE: several target values are digit strings that "look like" numbers
# (e.g. "-2.05") but must be stored as text (shared strings), exactly like
# they were before the edit. A plain .Value assignment would have Excel
# auto-convert such look-alike strings to real numbers, so those are
# entered with a leading apostrophe (the normal Excel way of forcing text
# entry for a numeric-looking value). Expressions that already contain a
# non-numeric character (x, y, ...) naturally stay text without it.

# --- Restricciones_del_lider ---
$ws = $wb.Worksheets.Item(2)
$ws.Range("A2").Value = "1.0499999999999998 - x"
$ws.Range("B2").Value = "'-2.05"
$ws.Range("D2").Value = "'0.24"
$ws.Range("A3").Value = "-1.05 + x"
$ws.Range("B3").Value = "'0.050000000000000044"
$ws.Range("D3").Value = "'0.72"

# --- Restricciones_del_follower ---
$ws = $wb.Worksheets.Item(3)
$ws.Range("A2").Value = "-1 + 0.3508771929824561y"
$ws.Range("B2").Value = "'0.0"
$ws.Range("D2").Value = "'0.47"
$ws.Range("E2").Value = "'3.1"
$ws.Range("F2").Value = "'0.6"
$ws.Range("A3").Value = "-6.661338147750939e-16 + 2.220446049250313e-16y"
$ws.Range("B3").Value = "'-0.9999999999999993"
$ws.Range("D3").Value = "'0.88"
$ws.Range("E3").Value = "'5.8999999999999995"
$ws.Range("F3").Value = "'5.699999999999999"

# --- Punto_modificado ---
$ws = $wb.Worksheets.Item(4)
$ws.Range("A2").Value = "'1.05"
$ws.Range("B2").Value = "'2.85"

# --- Vector_bf ---
$ws = $wb.Worksheets.Item(5)
$ws.Range("A2").Value = "'-1.9997872807017547"

# --- Vector_BF ---
$ws = $wb.Worksheets.Item(6)
$ws.Range("A2").Value = "'-0.6976749999999999"
$ws.Range("A3").Value = "'-1.8832192982456155"

# --- Vector_Alpha (this one is a genuine numeric cell, not a shared string) ---
$ws = $wb.Worksheets.Item(7)
$ws.Range("A2").Value = 1.71
